$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank "Points for grading" (column E) scores for rows 22 and 24,
# matching the full marks already recorded in column D for those rubric items.
$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 10

# Match the author's final viewport/selection state in the sheet view.
$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
